# "Prepay Loan" sheet gains a new automation step: a "waittopageload" row
# (value 4000) inserted right after the header row (before the existing
# "repaymenttransactiondate" row), pushing the rest of the steps down by one
# row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prepay Loan")
$ws.Activate()

$ws.Rows.Item(2).Insert() | Out-Null

$ws.Range("A2").Value = "waittopageload"
$ws.Range("B2").Value = 4000

$ws.Range("B5").Select() | Out-Null
